# Auto-generated edit script: updates cached market-price derived values
# across multiple worksheets, per scheduled runner refresh.
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")

# ALC!row48 (hunk 0)
$ws_ALC.Range("H48").Value = 0
$ws_ALC.Range("I48").Value = 0
$ws_ALC.Range("K48").Value = 0
$ws_ALC.Range("M48").ClearContents()

# ALC!row51 (hunk 1)
$ws_ALC.Range("H51").Value = 3564.4285
$ws_ALC.Range("J51").Value = 3858.5
$ws_ALC.Range("L51").Value = 3858.5
$ws_ALC.Range("N51").Value = -4826.5

# ALC!row56 (hunk 2)
$ws_ALC.Range("H56").Value = 0
$ws_ALC.Range("I56").Value = 0
$ws_ALC.Range("K56").Value = 0
$ws_ALC.Range("M56").ClearContents()

# ALC!row76 (hunk 3)
$ws_ALC.Range("H76").Value = 3475409.8
$ws_ALC.Range("I76").Value = 4118263.5
$ws_ALC.Range("J76").Value = 4000
$ws_ALC.Range("K76").Value = 4118263.5
$ws_ALC.Range("L76").Value = 4000
$ws_ALC.Range("M76").Value = -4117948.5
$ws_ALC.Range("N76").Value = -4630

# ALC!row79 (hunk 4)
$ws_ALC.Range("H79").Value = 3475409.8
$ws_ALC.Range("I79").Value = 4118263.5
$ws_ALC.Range("J79").Value = 4000
$ws_ALC.Range("K79").Value = 4118263.5
$ws_ALC.Range("L79").Value = 4000
$ws_ALC.Range("M79").Value = -4117171.5
$ws_ALC.Range("N79").Value = -6184

# ALC!row113 (hunk 5)
$ws_ALC.Range("H113").Value = 1988.5454
$ws_ALC.Range("I113").Value = 1984.25
$ws_ALC.Range("K113").Value = 1984.25
$ws_ALC.Range("M113").Value = 1269.75

# ALC!row116 (hunk 6)
$ws_ALC.Range("H116").Value = 9886672
$ws_ALC.Range("I116").Value = 13840271
$ws_ALC.Range("J116").Value = 2672.5
$ws_ALC.Range("K116").Value = 13840271
$ws_ALC.Range("L116").Value = 2672.5
$ws_ALC.Range("M116").Value = -13836829
$ws_ALC.Range("N116").Value = -9556.5

# ARM!row62 (hunk 7)
$ws_ARM.Range("H62").Value = 0
$ws_ARM.Range("J62").Value = 0
$ws_ARM.Range("L62").Value = 0
$ws_ARM.Range("N62").ClearContents()

# ARM!row63 (hunk 8)
$ws_ARM.Range("H63").Value = 6606.5
$ws_ARM.Range("I63").Value = 6606.5
$ws_ARM.Range("K63").Value = 6606.5
$ws_ARM.Range("M63").Value = -5920.5

# ARM!row65 (hunk 9)
$ws_ARM.Range("H65").Value = 0
$ws_ARM.Range("J65").Value = 0
$ws_ARM.Range("L65").Value = 0
$ws_ARM.Range("N65").ClearContents()

# ARM!row66 (hunk 10)
$ws_ARM.Range("H66").Value = 6606.5
$ws_ARM.Range("I66").Value = 6606.5
$ws_ARM.Range("K66").Value = 33032.5
$ws_ARM.Range("M66").Value = -29600.5

# ARM!row102 (hunk 11)
$ws_ARM.Range("H102").Value = 2046.6666
$ws_ARM.Range("I102").Value = 1816
$ws_ARM.Range("J102").Value = 3200
$ws_ARM.Range("K102").Value = 1816
$ws_ARM.Range("L102").Value = 3200
$ws_ARM.Range("M102").Value = -194
$ws_ARM.Range("N102").Value = -6444

# ARM!row110 (hunk 12)
$ws_ARM.Range("H110").Value = 1727.1111
$ws_ARM.Range("I110").Value = 772.5714
$ws_ARM.Range("J110").Value = 2334.5454
$ws_ARM.Range("K110").Value = 772.5714
$ws_ARM.Range("L110").Value = 2334.5454
$ws_ARM.Range("M110").Value = 1272.4286
$ws_ARM.Range("N110").Value = -6424.5454

# ARM!row134 (hunk 13)
$ws_ARM.Range("H134").Value = 35376.332
$ws_ARM.Range("J134").Value = 35376.332
$ws_ARM.Range("L134").Value = 35376.332
$ws_ARM.Range("N134").Value = -45516.332

# BSM!row105 (hunk 14)
$ws_BSM.Range("H105").Value = 347923.47
$ws_BSM.Range("I105").Value = 2963.5
$ws_BSM.Range("J105").Value = 1114501.2
$ws_BSM.Range("K105").Value = 2963.5
$ws_BSM.Range("L105").Value = 1114501.2
$ws_BSM.Range("M105").Value = -1216.5
$ws_BSM.Range("N105").Value = -1117995.2

# BSM!row107 (hunk 15)
$ws_BSM.Range("H107").Value = 982.9545000000001
$ws_BSM.Range("I107").Value = 705.3333
$ws_BSM.Range("J107").Value = 1577.8572
$ws_BSM.Range("K107").Value = 705.3333
$ws_BSM.Range("L107").Value = 1577.8572
$ws_BSM.Range("M107").Value = 1214.6667
$ws_BSM.Range("N107").Value = -5417.8572

# CRP!row95 (hunk 16)
$ws_CRP.Range("H95").Value = 9700
$ws_CRP.Range("J95").Value = 9700
$ws_CRP.Range("L95").Value = 9700
$ws_CRP.Range("N95").Value = -15192

# CRP!row99 (hunk 17)
$ws_CRP.Range("H99").Value = 5690318
$ws_CRP.Range("I99").Value = 6954500
$ws_CRP.Range("K99").Value = 6954500
$ws_CRP.Range("M99").Value = -6953002

# CRP!row105 (hunk 18)
$ws_CRP.Range("H105").Value = 862.8333
$ws_CRP.Range("I105").Value = 884.5454999999999
$ws_CRP.Range("J105").Value = 828.7143
$ws_CRP.Range("K105").Value = 884.5454999999999
$ws_CRP.Range("L105").Value = 828.7143
$ws_CRP.Range("M105").Value = 862.4545000000001
$ws_CRP.Range("N105").Value = -4322.7143

# CRP!row126 (hunk 19)
$ws_CRP.Range("H126").Value = 5690318
$ws_CRP.Range("I126").Value = 6954500
$ws_CRP.Range("K126").Value = 20863500
$ws_CRP.Range("M126").Value = -20861030

# CRP!row138 (hunk 20)
$ws_CRP.Range("H138").Value = 53804
$ws_CRP.Range("J138").Value = 53804
$ws_CRP.Range("L138").Value = 53804
$ws_CRP.Range("N138").Value = -64084

# CUL!row17 (hunk 21)
$ws_CUL.Range("H17").Value = 233.33333
$ws_CUL.Range("J17").Value = 300
$ws_CUL.Range("L17").Value = 900
$ws_CUL.Range("N17").Value = -1238

# CUL!row80 (hunk 22)
$ws_CUL.Range("H80").Value = 1165
$ws_CUL.Range("I80").Value = 0
$ws_CUL.Range("J80").Value = 1165
$ws_CUL.Range("K80").Value = 0
$ws_CUL.Range("L80").Value = 3495
$ws_CUL.Range("M80").ClearContents()
$ws_CUL.Range("N80").Value = -5367

# CUL!row83 (hunk 23)
$ws_CUL.Range("H83").Value = 1165
$ws_CUL.Range("I83").Value = 0
$ws_CUL.Range("J83").Value = 1165
$ws_CUL.Range("K83").Value = 0
$ws_CUL.Range("L83").Value = 10485
$ws_CUL.Range("M83").ClearContents()
$ws_CUL.Range("N83").Value = -19845

# CUL!row103 (hunk 24)
$ws_CUL.Range("H103").Value = 335832.22
$ws_CUL.Range("J103").Value = 503244.16
$ws_CUL.Range("L103").Value = 1509732.48
$ws_CUL.Range("N103").Value = -1511490.48

# CUL!row105 (hunk 25)
$ws_CUL.Range("H105").Value = 5666.6665
$ws_CUL.Range("J105").Value = 5666.6665
$ws_CUL.Range("L105").Value = 16999.9995
$ws_CUL.Range("N105").Value = -22241.9995

# GSM!row59 (hunk 26)
$ws_GSM.Range("H59").Value = 8496.666999999999

# GSM!row70 (hunk 27)
$ws_GSM.Range("H70").Value = 5521.3887
$ws_GSM.Range("I70").Value = 5488.4614
$ws_GSM.Range("J70").Value = 5607
$ws_GSM.Range("K70").Value = 5488.4614
$ws_GSM.Range("L70").Value = 5607
$ws_GSM.Range("M70").Value = -5218.4614
$ws_GSM.Range("N70").Value = -6147

# GSM!row73 (hunk 28)
$ws_GSM.Range("H73").Value = 5521.3887
$ws_GSM.Range("I73").Value = 5488.4614
$ws_GSM.Range("J73").Value = 5607
$ws_GSM.Range("K73").Value = 5488.4614
$ws_GSM.Range("L73").Value = 5607
$ws_GSM.Range("M73").Value = -4552.4614
$ws_GSM.Range("N73").Value = -7479

# GSM!row97 (hunk 29)
$ws_GSM.Range("H97").Value = 989
$ws_GSM.Range("I97").Value = 989
$ws_GSM.Range("J97").Value = 0
$ws_GSM.Range("K97").Value = 989
$ws_GSM.Range("L97").Value = 0
$ws_GSM.Range("M97").Value = -493
$ws_GSM.Range("N97").ClearContents()

# GSM!row107 (hunk 30)
$ws_GSM.Range("H107").Value = 199.27777
$ws_GSM.Range("I107").Value = 152.17647
$ws_GSM.Range("J107").Value = 1000
$ws_GSM.Range("K107").Value = 152.17647
$ws_GSM.Range("L107").Value = 1000
$ws_GSM.Range("M107").Value = 1767.82353
$ws_GSM.Range("N107").Value = -4840

# GSM!row138 (hunk 31)
$ws_GSM.Range("H138").Value = 66304.75
$ws_GSM.Range("J138").Value = 66304.75
$ws_GSM.Range("L138").Value = 66304.75
$ws_GSM.Range("N138").Value = -76584.75

# LTW!row16 (hunk 32)
$ws_LTW.Range("H16").Value = 4167994
$ws_LTW.Range("I16").Value = 7693442
$ws_LTW.Range("K16").Value = 7693442
$ws_LTW.Range("M16").Value = -7693272

# LTW!row61 (hunk 33)
$ws_LTW.Range("H61").Value = 2391.5
$ws_LTW.Range("I61").Value = 2519.8
$ws_LTW.Range("J61").Value = 1750
$ws_LTW.Range("K61").Value = 2519.8
$ws_LTW.Range("L61").Value = 1750
$ws_LTW.Range("M61").Value = -2317.8
$ws_LTW.Range("N61").Value = -2154

# LTW!row82 (hunk 34)
$ws_LTW.Range("H82").Value = 1651
$ws_LTW.Range("I82").Value = 1602
$ws_LTW.Range("J82").Value = 1667.3334
$ws_LTW.Range("K82").Value = 1602
$ws_LTW.Range("L82").Value = 1667.3334
$ws_LTW.Range("M82").Value = -1241
$ws_LTW.Range("N82").Value = -2389.3334

# LTW!row85 (hunk 35)
$ws_LTW.Range("H85").Value = 1651
$ws_LTW.Range("I85").Value = 1602
$ws_LTW.Range("J85").Value = 1667.3334
$ws_LTW.Range("K85").Value = 1602
$ws_LTW.Range("L85").Value = 1667.3334
$ws_LTW.Range("M85").Value = -354
$ws_LTW.Range("N85").Value = -4163.3334

# LTW!row93 (hunk 36)
$ws_LTW.Range("H93").Value = 1611.3043
$ws_LTW.Range("I93").Value = 1626.6666
$ws_LTW.Range("J93").Value = 1450
$ws_LTW.Range("K93").Value = 1626.6666
$ws_LTW.Range("L93").Value = 1450
$ws_LTW.Range("M93").Value = -378.6666
$ws_LTW.Range("N93").Value = -3946

# LTW!row106 (hunk 37)
$ws_LTW.Range("H106").Value = 21881.75
$ws_LTW.Range("J106").Value = 21881.75
$ws_LTW.Range("L106").Value = 21881.75
$ws_LTW.Range("N106").Value = -24405.75

# LTW!row113 (hunk 38)
$ws_LTW.Range("H113").Value = 2391.5
$ws_LTW.Range("I113").Value = 2519.8
$ws_LTW.Range("J113").Value = 1750
$ws_LTW.Range("K113").Value = 2519.8
$ws_LTW.Range("L113").Value = 1750
$ws_LTW.Range("M113").Value = -349.8000000000002
$ws_LTW.Range("N113").Value = -6090

# LTW!row123 (hunk 39)
$ws_LTW.Range("H123").Value = 41000
$ws_LTW.Range("J123").Value = 41000
$ws_LTW.Range("L123").Value = 41000
$ws_LTW.Range("N123").Value = -50800
